$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns hold plain text in the source data (e.g.
# "58.867.24" or "  +2.40%  "). Some of the updated Price strings look
# like ordinary decimals (e.g. "139.60"), so a direct .Value assignment
# would be auto-coerced to a number and silently lose the trailing zero
# / original text formatting. For those cells we briefly force a text
# number format, assign the literal string, then restore the default
# "Normal" style so no stray formatting is left behind.

$ws.Range("D2").Value = "58.781.23"
$ws.Range("E2").Value = "  +2.28%  "
$ws.Range("D3").Value = "2.585.49"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.81%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("D9").Value = "2.595.54"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.100"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").Value = "3.042.90"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").Value = "58.777.67"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.73%  "
$ws.Range("D17").Value = "2.582.92"
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "0.0₃0721"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("E31").Value = "  -4.57%  "
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.825"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.815"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "273.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.589"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0520"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "1.988.45"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0221"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("E51").Value = "  -3.16%  "
